$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(25).Insert()

$ws.Range("A25").Value = 11
$ws.Range("B25").Value = "Vega Monumental Concepción"
$ws.Range("C25").Value = "Bíobío"
$ws.Range("D25").Value = 45260
$ws.Range("E25").Value = 8
$ws.Range("F25").Value = 100112031
$ws.Range("G25").Value = "Poroto verde"
$ws.Range("H25").Value = "Magnum"
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 100
$ws.Range("K25").Value = 39000
$ws.Range("L25").Value = 39000
$ws.Range("M25").Value = 39000
$ws.Range("N25").Value = "$/malla 25 kilos"
$ws.Range("O25").Value = "Región de Coquimbo"
$ws.Range("P25").Value = 1560
$ws.Range("Q25").Value = 25
$ws.Range("R25").Value = "Hortaliza"

$ws.Range("D25").NumberFormat = "YYYY-MM-DD HH:MM:SS"
